# Natmi following Dr Hou advice
# Update computed L-R interaction values in rows 2-7 (columns E..T)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; F=1; G=1.141592333333333;  H=3.424777;           I=0.5720393989604073; J=0.5720393989604073; K=3; L=1; M=122.298789;         N=366.8963669999999;  O=0.4294659933718775; P=0.4294659933718776; Q=139.615359898351;  R=1256.538239085159; S=0.2456714687223831; T=0.2456714687223831 }
    3 = @{ E=3; F=1; G=1.141592333333333;  H=3.424777;           I=0.5720393989604073; J=0.5720393989604073; K=3; L=1; M=119.3038053333333;  N=357.911416;         O=0.4189487703800985; P=0.4189487703800985; Q=136.1963095060258; R=1225.766785554232; S=0.2396552028034332; T=0.2396552028034333 }
    4 = @{ E=3; F=1; G=1.141592333333333;  H=3.424777;           I=0.5720393989604073; J=0.5720393989604073; K=3; L=1; M=43.16684233333334;  N=129.500527;         O=0.1515852362480238; P=0.1515852362480239; Q=49.27893626194211; R=443.510426357479;  S=0.08671272743459091; T=0.08671272743459092 }
    5 = @{ E=3; F=1; G=0.8540610000000001; H=2.562183;           I=0.4279606010395928; J=0.4279606010395928; K=3; L=1; M=122.298789;         N=366.8963669999999;  O=0.4294659933718775; P=0.4294659933718776; Q=104.450626032129;  R=940.0556342891609; S=0.1837945246494945; T=0.1837945246494945 }
    6 = @{ E=3; F=1; G=0.8540610000000001; H=2.562183;           I=0.4279606010395928; J=0.4279606010395928; K=3; L=1; M=119.3038053333333;  N=357.911416;         O=0.4189487703800985; P=0.4189487703800985; Q=101.892727286792;  R=917.0345455811281; S=0.1792935675766653; T=0.1792935675766653 }
    7 = @{ E=3; F=1; G=0.8540610000000001; H=2.562183;           I=0.4279606010395928; J=0.4279606010395928; K=3; L=1; M=43.16684233333334;  N=129.500527;         O=0.1515852362480238; P=0.1515852362480239; Q=36.867116530049;   R=331.804048770441;  S=0.06487250881343294; T=0.06487250881343296 }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
